# This edit adds a new weekly price record for "Haba" (Femacal de La Calera)
# to the sheet. The new record is inserted as row 94, which pushes all the
# existing records that used to be rows 94-111 down by one row (to 95-112),
# preserving their original data/order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting rows 94:111 down to 95:112.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly record.
$ws.Cells.Item(94, 1).Value = 3
$ws.Cells.Item(94, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(94, 3).Value = "Coquimbo"
$ws.Cells.Item(94, 4).Value = 44543
$ws.Cells.Item(94, 5).Value = 5
$ws.Cells.Item(94, 6).Value = 100112026
$ws.Cells.Item(94, 7).Value = "Haba"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 100
$ws.Cells.Item(94, 11).Value = 8500
$ws.Cells.Item(94, 12).Value = 9000
$ws.Cells.Item(94, 13).Value = 8700
$ws.Cells.Item(94, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(94, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(94, 16).Value = 348
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = "Hortaliza"
